# Applies the "Made Screenshot, Created GeneralGameDocument" edit:
#  - re-splits several run texts and wraps the "interesting" (foreign /
#    proper-noun) words in <w:proofErr spellStart/spellEnd> pairs, as Word's
#    background spell checker does after you retype / re-enter a paragraph
#  - appends " (schneller stärker)" to the "Boss Gegner" bullet
#  - moves the "_GoBack" last-edit bookmark from the "Particle effect"
#    bullet to the end of the (now longer) "Boss Gegner" bullet
#
# Strategy: Word's COM object model has no property for <w:proofErr/>, so we
# drive the change at the OOXML level through Range.InsertXML, which (like
# real Word) replaces only the content of the given range and leaves the
# paragraph's own <w:pPr> alone as long as the range excludes the trailing
# paragraph mark.

$d = $word.ActiveDocument

function Set-ParagraphRuns($Index, $InnerXml) {
    $para = $d.Paragraphs($Index)
    $full = $para.Range
    $target = $d.Range($full.Start, $full.End - 1)
    $pkg = "<?xml version='1.0' encoding='UTF-8' standalone='yes'?>" +
        "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
        "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
        "<pkg:xmlData>" +
        "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
        "<w:body><w:p>$InnerXml</w:p></w:body>" +
        "</w:document>" +
        "</pkg:xmlData></pkg:part></pkg:package>"
    $target.InsertXML($pkg)
}

# 1) Remove the stray "_GoBack" bookmark currently sitting after "Particle
#    effect" - it is being relocated onto "Boss Gegner" below.
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
}

# 2) "Items/power ups" -> "Items/power " + spell-checked "ups"
$xml2 = '<w:r><w:t xml:space="preserve">Items/power </w:t></w:r>'
$xml2 = $xml2 + '<w:proofErr w:type="spellStart"/>'
$xml2 = $xml2 + '<w:r><w:t>ups</w:t></w:r>'
$xml2 = $xml2 + '<w:proofErr w:type="spellEnd"/>'
Set-ParagraphRuns 3 $xml2

# 3) "Boss Gegner" -> "Boss Gegner (schneller stärker)" and gains the
#    relocated "_GoBack" bookmark.
$xml3 = '<w:r><w:t>Boss Gegner</w:t></w:r>'
$xml3 = $xml3 + '<w:r><w:t xml:space="preserve"> (schneller stärker)</w:t></w:r>'
$xml3 = $xml3 + '<w:bookmarkStart w:id="0" w:name="_GoBack"/>'
$xml3 = $xml3 + '<w:bookmarkEnd w:id="0"/>'
Set-ParagraphRuns 5 $xml3

# 4) "Cutscene Einstieg, Roboter spawnen, ..." -> spell-checked "Cutscene"
#    and "spawnen".
$xml4 = '<w:proofErr w:type="spellStart"/>'
$xml4 = $xml4 + '<w:r><w:t>Cutscene</w:t></w:r>'
$xml4 = $xml4 + '<w:proofErr w:type="spellEnd"/>'
$xml4 = $xml4 + '<w:r><w:t xml:space="preserve"> Einstieg, Roboter </w:t></w:r>'
$xml4 = $xml4 + '<w:proofErr w:type="spellStart"/>'
$xml4 = $xml4 + '<w:r><w:t>spawnen</w:t></w:r>'
$xml4 = $xml4 + '<w:proofErr w:type="spellEnd"/>'
$xml4 = $xml4 + '<w:r><w:t>, Türen öffnen sich. Spieler bekommt gesagt: „Lass sie nicht entkommen!“</w:t></w:r>'
Set-ParagraphRuns 7 $xml4

# 5) "Platforms move" -> spell-checked "Platforms" and "move".
$xml5 = '<w:proofErr w:type="spellStart"/>'
$xml5 = $xml5 + '<w:r><w:t>Platforms</w:t></w:r>'
$xml5 = $xml5 + '<w:proofErr w:type="spellEnd"/>'
$xml5 = $xml5 + '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
$xml5 = $xml5 + '<w:proofErr w:type="spellStart"/>'
$xml5 = $xml5 + '<w:r><w:t>move</w:t></w:r>'
$xml5 = $xml5 + '<w:proofErr w:type="spellEnd"/>'
Set-ParagraphRuns 8 $xml5

# 6) "Spawn points move" -> spell-checked "Spawn", "points" and "move".
$xml6 = '<w:proofErr w:type="spellStart"/>'
$xml6 = $xml6 + '<w:r><w:t>Spawn</w:t></w:r>'
$xml6 = $xml6 + '<w:proofErr w:type="spellEnd"/>'
$xml6 = $xml6 + '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
$xml6 = $xml6 + '<w:proofErr w:type="spellStart"/>'
$xml6 = $xml6 + '<w:r><w:t>points</w:t></w:r>'
$xml6 = $xml6 + '<w:proofErr w:type="spellEnd"/>'
$xml6 = $xml6 + '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
$xml6 = $xml6 + '<w:proofErr w:type="spellStart"/>'
$xml6 = $xml6 + '<w:r><w:t>move</w:t></w:r>'
$xml6 = $xml6 + '<w:proofErr w:type="spellEnd"/>'
Set-ParagraphRuns 9 $xml6

# 7) "randomize delay between impact sfx and enemy death sfx" -> both "sfx"
#    occurrences spell-checked; every run keeps the bold/bCs/en-US rPr.
$rPr = '<w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr>'
$xml7 = "<w:r>$rPr" + '<w:t xml:space="preserve">randomize delay between impact </w:t></w:r>'
$xml7 = $xml7 + '<w:proofErr w:type="spellStart"/>'
$xml7 = $xml7 + "<w:r>$rPr" + '<w:t>sfx</w:t></w:r>'
$xml7 = $xml7 + '<w:proofErr w:type="spellEnd"/>'
$xml7 = $xml7 + "<w:r>$rPr" + '<w:t xml:space="preserve"> and enemy death </w:t></w:r>'
$xml7 = $xml7 + '<w:proofErr w:type="spellStart"/>'
$xml7 = $xml7 + "<w:r>$rPr" + '<w:t>sfx</w:t></w:r>'
$xml7 = $xml7 + '<w:proofErr w:type="spellEnd"/>'
Set-ParagraphRuns 16 $xml7
